$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships protected; unprotect so the model-holdings cells
# (which live on a protected sheet) can be updated, then re-protect at
# the end to restore the original protected state.
$ws.Unprotect()

# Footer disclaimer text: bump the "as of" model date from 2021-03-24 to
# 2021-03-25 (confidential client disclosure string in A11).
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

# Refreshed Weight (D) / Percent Change (E) figures for each holding row.
$ws.Range("D2").Value = 0.4997971018581381
$ws.Range("E2").Value = 0.007702829169966119

$ws.Range("D3").Value = 0.2448572731675436
$ws.Range("E3").Value = -0.002197457228064637

$ws.Range("D4").Value = 0.09722490538404785
$ws.Range("E4").Value = 0.01443899854285347

$ws.Range("D5").Value = 0.1007846689550715
$ws.Range("E5").Value = 0.02160243407707907

$ws.Range("D6").Value = 0.02970796773336915
$ws.Range("E6").Value = 0.02652766308835686

$ws.Range("D7").Value = 0.02762808290182981
$ws.Range("E7").Value = 0.020836724727332

$ws.Range("E8").Value = 0.008256574461982691

# Restore sheet protection (original password hash token was "D382").
$ws.Protect("D382")
